# Update Fbn1-Itgav NATMI ligand-receptor output with recomputed TPM-based
# expression values. Only the "ECs" cluster's raw ligand (Fbn1) and receptor
# (Itgav) expression figures changed; every other column here (specificity
# scores, edge weights, edge specificities) is a downstream recalculation of
# those raw values across all Sending/Target cluster combinations, so each
# affected cell is written explicitly with its new value below.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 2.763564666666667
$ws.Range("H2").Value = 8.290694
$ws.Range("I2").Value = 0.009909756414635561
$ws.Range("J2").Value = 0.009909756414635559
$ws.Range("M2").Value = 3.759736666666667
$ws.Range("N2").Value = 11.27921
$ws.Range("O2").Value = 0.0683751702595819
$ws.Range("P2").Value = 0.06837517025958188
$ws.Range("Q2").Value = 10.39027540797111
$ws.Range("R2").Value = 93.51247867174
$ws.Range("S2").Value = 0.0006775812820816903
$ws.Range("T2").Value = 0.0006775812820816901
$ws.Range("G3").Value = 2.763564666666667
$ws.Range("H3").Value = 8.290694
$ws.Range("I3").Value = 0.009909756414635561
$ws.Range("J3").Value = 0.009909756414635559
$ws.Range("O3").Value = 0.6514180024294648
$ws.Range("P3").Value = 0.6514180024294647
$ws.Range("Q3").Value = 98.98933231546911
$ws.Range("R3").Value = 890.903990839222
$ws.Range("S3").Value = 0.006455393728184472
$ws.Range("T3").Value = 0.00645539372818447
$ws.Range("G4").Value = 2.763564666666667
$ws.Range("H4").Value = 8.290694
$ws.Range("I4").Value = 0.009909756414635561
$ws.Range("J4").Value = 0.009909756414635559
$ws.Range("O4").Value = 0.2802068273109533
$ws.Range("P4").Value = 0.2802068273109533
$ws.Range("Q4").Value = 42.58016610271778
$ws.Range("R4").Value = 383.22149492446
$ws.Range("S4").Value = 0.002776781404369399
$ws.Range("T4").Value = 0.002776781404369397
$ws.Range("G5").Value = 266.1315866666666
$ws.Range("H5").Value = 798.3947599999999
$ws.Range("I5").Value = 0.9543106517164206
$ws.Range("J5").Value = 0.9543106517164204
$ws.Range("M5").Value = 3.759736666666667
$ws.Range("N5").Value = 11.27921
$ws.Range("O5").Value = 0.0683751702595819
$ws.Range("P5").Value = 0.06837517025958188
$ws.Range("Q5").Value = 1000.584684548844
$ws.Range("R5").Value = 9005.262160939597
$ws.Range("S5").Value = 0.06525115329164281
$ws.Range("T5").Value = 0.06525115329164279
$ws.Range("G6").Value = 266.1315866666666
$ws.Range("H6").Value = 798.3947599999999
$ws.Range("I6").Value = 0.9543106517164206
$ws.Range("J6").Value = 0.9543106517164204
$ws.Range("O6").Value = 0.6514180024294648
$ws.Range("P6").Value = 0.6514180024294647
$ws.Range("Q6").Value = 9532.683779737763
$ws.Range("R6").Value = 85794.15401763987
$ws.Range("S6").Value = 0.6216551384382714
$ws.Range("T6").Value = 0.6216551384382711
$ws.Range("G7").Value = 266.1315866666666
$ws.Range("H7").Value = 798.3947599999999
$ws.Range("I7").Value = 0.9543106517164206
$ws.Range("J7").Value = 0.9543106517164204
$ws.Range("O7").Value = 0.2802068273109533
$ws.Range("P7").Value = 0.2802068273109533
$ws.Range("Q7").Value = 4100.47476077871
$ws.Range("S7").Value = 0.2674043599865064
$ws.Range("T7").Value = 0.2674043599865062
$ws.Range("I8").Value = 0.03577959186894402
$ws.Range("J8").Value = 0.03577959186894401
$ws.Range("M8").Value = 3.759736666666667
$ws.Range("N8").Value = 11.27921
$ws.Range("O8").Value = 0.0683751702595819
$ws.Range("P8").Value = 0.06837517025958188
$ws.Range("Q8").Value = 37.51452588219888
$ws.Range("R8").Value = 337.63073293979
$ws.Range("S8").Value = 0.002446435685857399
$ws.Range("T8").Value = 0.002446435685857398
$ws.Range("I9").Value = 0.03577959186894402
$ws.Range("J9").Value = 0.03577959186894401
$ws.Range("O9").Value = 0.6514180024294648
$ws.Range("P9").Value = 0.6514180024294647
$ws.Range("S9").Value = 0.02330747026300903
$ws.Range("T9").Value = 0.02330747026300902
$ws.Range("I10").Value = 0.03577959186894402
$ws.Range("J10").Value = 0.03577959186894401
$ws.Range("O10").Value = 0.2802068273109533
$ws.Range("P10").Value = 0.2802068273109533
$ws.Range("S10").Value = 0.01002568592007759
$ws.Range("T10").Value = 0.01002568592007758
